$wb = $excel.ActiveWorkbook

# --- Rename existing headers ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy the bold/centered header style from an existing sheet's header row
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# --- Data rows (A2:D31) ---
$data = New-Object 'object[,]' 30,4
$data[0,0] = 45039.99999999999
$data[0,1] = 27
$data[0,2] = -202.8635727769995
$data[0,3] = 284.1561925255714
$data[1,0] = 45053.99999999999
$data[1,1] = 47
$data[1,2] = -192.8687157905289
$data[1,3] = 287.5140429485451
$data[2,0] = 45060.99999999999
$data[2,1] = 57
$data[2,2] = -189.7453598415397
$data[2,3] = 291.8990598294594
$data[3,0] = 45067.99999999999
$data[3,1] = 66
$data[3,2] = -178.0588071617111
$data[3,3] = 312.9125579642754
$data[4,0] = 45074.99999999999
$data[4,1] = 76
$data[4,2] = -176.3425100633189
$data[4,3] = 311.3081505737908
$data[5,0] = 45081.99999999999
$data[5,1] = 86
$data[5,2] = -148.6624365572835
$data[5,3] = 332.4097344432367
$data[6,0] = 45088.99999999999
$data[6,1] = 96
$data[6,2] = -156.1322014449593
$data[6,3] = 334.2309246189285
$data[7,0] = 45095.99999999999
$data[7,1] = 106
$data[7,2] = -132.7747125672468
$data[7,3] = 334.0601070366855
$data[8,0] = 45102.99999999999
$data[8,1] = 116
$data[8,2] = -119.6655386299419
$data[8,3] = 371.4929316892905
$data[9,0] = 45123.99999999999
$data[9,1] = 145
$data[9,2] = -108.7835160166261
$data[9,3] = 394.1344139254093
$data[10,0] = 45130.99999999999
$data[10,1] = 155
$data[10,2] = -88.10787608930444
$data[10,3] = 408.4624461893928
$data[11,0] = 45144.99999999999
$data[11,1] = 175
$data[11,2] = -75.5032081500615
$data[11,3] = 414.2430759986129
$data[12,0] = 45151.99999999999
$data[12,1] = 185
$data[12,2] = -62.87848337455383
$data[12,3] = 418.8418665946095
$data[13,0] = 45172.99999999999
$data[13,1] = 214
$data[13,2] = -25.09451950646735
$data[13,3] = 471.689708558026
$data[14,0] = 45179.99999999999
$data[14,1] = 224
$data[14,2] = -4.889856550039705
$data[14,3] = 467.106979003664
$data[15,0] = 45186.99999999999
$data[15,1] = 234
$data[15,2] = 8.173508247153569
$data[15,3] = 476.9547651488971
$data[16,0] = 45193.99999999999
$data[16,1] = 244
$data[16,2] = 32.9041982393868
$data[16,3] = 482.7621548996375
$data[17,0] = 45200.99999999999
$data[17,1] = 254
$data[17,2] = 21.78281918669733
$data[17,3] = 481.9480631763515
$data[18,0] = 45228.99999999999
$data[18,1] = 293
$data[18,2] = 48.68833810964768
$data[18,3] = 546.3291677403976
$data[19,0] = 45235.99999999999
$data[19,1] = 303
$data[19,2] = 60.45622478727026
$data[19,3] = 544.3824832138426
$data[20,0] = 45242.99999999999
$data[20,1] = 313
$data[20,2] = 58.96182542808541
$data[20,3] = 543.3745447467019
$data[21,0] = 45249.99999999999
$data[21,1] = 323
$data[21,2] = 60.86403329918749
$data[21,3] = 558.7319112537994
$data[22,0] = 45256.99999999999
$data[22,1] = 333
$data[22,2] = 83.44239008115072
$data[22,3] = 582.5509353217151
$data[23,0] = 45263.99999999999
$data[23,1] = 343
$data[23,2] = 92.83520857477819
$data[23,3] = 591.829623457608
$data[24,0] = 45270.99999999999
$data[24,1] = 353
$data[24,2] = 116.4934122145539
$data[24,3] = 588.7160980185839
$data[25,0] = 45277.99999999999
$data[25,1] = 363
$data[25,2] = 128.8234941792677
$data[25,3] = 606.1585475680794
$data[26,0] = 45284.99999999999
$data[26,1] = 372
$data[26,2] = 126.4628130631263
$data[26,3] = 606.9301286599317
$data[27,0] = 45291.99999999999
$data[27,1] = 382
$data[27,2] = 137.5123834220825
$data[27,3] = 613.0109405071452
$data[28,0] = 45298.99999999999
$data[28,1] = 392
$data[28,2] = 144.4752250248308
$data[28,3] = 618.731045328993
$data[29,0] = 45305.99999999999
$data[29,1] = 402
$data[29,2] = 159.3371444400644
$data[29,3] = 645.2468979219401
$wsForecast.Range("A2:D31").Value = $data

# Copy the date number-format style from an existing sheet's date column
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A31").PasteSpecial(-4122)

$excel.CutCopyMode = 0

Write-Output "Done applying PO Forecast sheet changes"
